$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff shows cell C10 changing from an integer value 18 to a
# numeric value of 1 (explicitly typed as numeric).
$ws.Range("C10").Value = 1
